# Updates res_bus/vm_pu.xlsx results for the 380 kV case (commit: "case with 380 kV done").
# Column B changes from the flat 1.05 pu start value to 1.02 pu, and the dependent
# bus-voltage magnitudes (C-F, I-N) are replaced with the recomputed power-flow results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$newValues = [ordered]@{
    2 = @(1.02, 1.052341433164792, 1.06233296558097, 1.059486112305843, 1.070943064865677, 1.043765718785753, 1.057364545090457, 1.065054597941509, 1.062215485499374, 1.073641556505679, 1.058866124025895)
    3 = @(1.02, 1.053383442592033, 1.063305433459447, 1.06041129721509, 1.071949977632898, 1.043967765682469, 1.058056549299484, 1.065841464785452, 1.062954628572978, 1.074464457026014, 1.059559110960272)
    4 = @(1.02, 1.054057755711031, 1.063935012681823, 1.061010352412002, 1.072601897537258, 1.04409659364622, 1.058503780437783, 1.066350326414741, 1.063432671491763, 1.074996692514296, 1.060006977218085)
    5 = @(1.02, 1.05434125184186, 1.064199765453751, 1.061262289887336, 1.072876054879651, 1.044150295211857, 1.058691666134038, 1.066564180876642, 1.063633584562795, 1.075220387098764, 1.060195129733579)
    6 = @(1.02, 1.054388852955332, 1.064244223142196, 1.061304596826865, 1.07292209237181, 1.044159285086653, 1.058723205313698, 1.066600083802093, 1.063667315439643, 1.075257943029017, 1.060226713702493)
    7 = @(1.02, 1.054061543742015, 1.063938550017069, 1.061013718441873, 1.072605560487568, 1.044097313008266, 1.058506291488243, 1.066353184228171, 1.063435356323548, 1.074999681757953, 1.060009491834526)
    8 = @(1.02, 1.052693571853296, 1.062661547408538, 1.059798699858941, 1.071283276609161, 1.043834396649203, 1.057598523017506, 1.065320583743374, 1.06246533004739, 1.07391970808701, 1.059100434228437)
    9 = @(1.02, 1.050283526967852, 1.060413845503032, 1.057660769240122, 1.068956186967405, 1.043356497710454, 1.055994792875045, 1.063498786270538, 1.06075426895982, 1.072014876897446, 1.057494426605224)
    10 = @(1.02, 1.04867717612191, 1.058917125598136, 1.056237601346218, 1.067406814968311, 1.043028104402945, 1.054922895352912, 1.06228279375918, 1.059612418165277, 1.070743830043101, 1.056421006865661)
    11 = @(1.02, 1.047981691999351, 1.058269451704024, 1.05562186503591, 1.066736407237528, 1.042883588687209, 1.054458106691662, 1.061755915550166, 1.059117719431313, 1.070193184822947, 1.055955558151199)
    12 = @(1.02, 1.047723369710183, 1.058028939773416, 1.055393229596774, 1.066487460718673, 1.042829561030624, 1.054285366211275, 1.061560158063069, 1.058933926202884, 1.069988610058496, 1.055782572359517)
    13 = @(1.02, 1.047778780198565, 1.058080527539904, 1.055442269205757, 1.066540857290413, 1.042841165892291, 1.054322424026747, 1.061602151030337, 1.058973352281467, 1.070032493869968, 1.055819682801338)
    14 = @(1.02, 1.047960338764135, 1.058249569624717, 1.055602964396308, 1.066715827753056, 1.042879129843054, 1.054443829891945, 1.0617397352225, 1.059102527841036, 1.070176275434788, 1.05594126107679)
    15 = @(1.02, 1.048072204509905, 1.058353730395218, 1.055701984154166, 1.066823642488296, 1.042902474571979, 1.054518619173337, 1.06182449857802, 1.059182111867203, 1.070264858635182, 1.056016156567539)
    16 = @(1.02, 1.048723334722432, 1.058960118385706, 1.056278476424097, 1.067451317894596, 1.043037646601621, 1.054953728183564, 1.062317753690195, 1.059645244044419, 1.070780368831066, 1.056451883482469)
    17 = @(1.02, 1.049131791957873, 1.059340601227265, 1.056640230366538, 1.06784517115447, 1.043121815768804, 1.055226486821601, 1.062627067418693, 1.059935682703067, 1.071103661865489, 1.056725029469057)
    18 = @(1.02, 1.049370045523596, 1.059562570504926, 1.05685128398436, 1.068074945428014, 1.043170686422755, 1.055385519583941, 1.062807451525676, 1.060105064597098, 1.071292206797456, 1.056884288076162)
    19 = @(1.02, 1.049451285082579, 1.059638263073151, 1.056923256049638, 1.068153300358676, 1.043187312082049, 1.05543973499353, 1.062868952173255, 1.060162815011682, 1.071356491249804, 1.056938580477851)
    20 = @(1.02, 1.049087967599556, 1.059299774870254, 1.056601412554154, 1.067802909630725, 1.043112808359509, 1.055197228887452, 1.062593884430039, 1.059904524076076, 1.071068978312205, 1.056695729985286)
    21 = @(1.02, 1.047906873974871, 1.058199789200844, 1.0555556415529, 1.066664301299745, 1.042867960011803, 1.054408081572953, 1.061699221524203, 1.059064489961562, 1.070133936464834, 1.055905461991084)
    22 = @(1.02, 1.047164339958031, 1.05750854926342, 1.054898566385035, 1.065948834306028, 1.042712000396859, 1.053911350702695, 1.061136413921368, 1.058536095428515, 1.069545803810877, 1.055408025706014)
    23 = @(1.02, 1.047557964992921, 1.057874953884423, 1.055246852212763, 1.066328076779827, 1.042794868301429, 1.054174730343259, 1.061434796937911, 1.058816229187865, 1.069857606147848, 1.055671779375878)
    24 = @(1.02, 1.049107769917755, 1.059318222419768, 1.056618952502287, 1.067822005650978, 1.043116879111528, 1.055210449482288, 1.062608878497281, 1.059918603398887, 1.071084650398597, 1.056708969354883)
    25 = @(1.02, 1.050906521268724, 1.060994625157143, 1.058213105189737, 1.069557442247121, 1.043481774240794, 1.056409881406954, 1.063970024549123, 1.061196823252413, 1.07250752846651, 1.057910104610458)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}
